$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 - this shifts existing rows 18..90 down to 19..91
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with its data
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C18").Value = "Los Lagos"
$ws.Range("D18").Value = 44575
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 100112052
$ws.Range("G18").Value = "Albahaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("N18").Value = "$/docena de matas"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 1167
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = "Hortaliza"
